$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting rows 11:99 down to 12:100
$ws.Rows.Item(11).Insert()

# Fill in the new row 11 with the same constant columns as all other data rows,
# plus the new record's specific values.
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44530
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112040
$ws.Cells.Item(11, 7).Value = "Cilantro"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 3300
$ws.Cells.Item(11, 11).Value = 1500
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 1750
$ws.Cells.Item(11, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 1167
$ws.Cells.Item(11, 17).Value = 1.5
$ws.Cells.Item(11, 18).Value = "Hortaliza"
